$wb = $excel.ActiveWorkbook

# --- Sheet "scripts汉化状态": rebuild the summary table (revert to the pre-rename,
# alphabetically-sorted layout with plain "1" flags instead of "Yes" text) ---
$ws3 = $wb.Worksheets.Item(3)

# Drop the old AutoFilter before touching the range beneath it.
$ws3.AutoFilterMode = $false

# Clear the old data block (rows 2-32, columns A-G) completely -- contents AND
# formatting -- so stale "Yes" shared-string cells / center-align styling from
# columns that no longer carry a mark in a given row are gone.
$ws3.Range("A2:G32").Clear()

# Row 2: Adv-Solar-Panel
$ws3.Cells.Item(2, 1).Value = "Adv-Solar-Panel"
$ws3.Cells.Item(2, 2).Value = 1
$ws3.Cells.Item(2, 3).Value = 1
$ws3.Cells.Item(2, 4).Value = 1
$ws3.Cells.Item(2, 4).HorizontalAlignment = -4108

# Row 3: Advanced-Backpacks
$ws3.Cells.Item(3, 1).Value = "Advanced-Backpacks"
$ws3.Cells.Item(3, 2).Value = 1
$ws3.Cells.Item(3, 3).Value = 6
$ws3.Cells.Item(3, 5).Value = 1
$ws3.Cells.Item(3, 5).HorizontalAlignment = -4108

# Row 4: AE2
$ws3.Cells.Item(4, 1).Value = "AE2"
$ws3.Cells.Item(4, 2).Value = 1
$ws3.Cells.Item(4, 3).Value = 1
$ws3.Cells.Item(4, 4).Value = 1
$ws3.Cells.Item(4, 4).HorizontalAlignment = -4108

# Row 5: Blood-Magic-Thaumcraft
$ws3.Cells.Item(5, 1).Value = "Blood-Magic-Thaumcraft"
$ws3.Cells.Item(5, 2).Value = 1
$ws3.Cells.Item(5, 3).Value = 870
$ws3.Cells.Item(5, 6).Value = 1
$ws3.Cells.Item(5, 6).HorizontalAlignment = -4108
$ws3.Cells.Item(5, 7).Value = 1
$ws3.Cells.Item(5, 7).HorizontalAlignment = -4108

# Row 6: Build-Craft
$ws3.Cells.Item(6, 1).Value = "Build-Craft"
$ws3.Cells.Item(6, 2).Value = 1
$ws3.Cells.Item(6, 3).Value = 15
$ws3.Cells.Item(6, 4).Value = 1
$ws3.Cells.Item(6, 4).HorizontalAlignment = -4108

# Row 7: CoreMod
$ws3.Cells.Item(7, 1).Value = "CoreMod"
$ws3.Cells.Item(7, 2).Value = 1
$ws3.Cells.Item(7, 3).Value = 326
$ws3.Cells.Item(7, 7).Value = 1
$ws3.Cells.Item(7, 7).HorizontalAlignment = -4108

# Row 8: Emt
$ws3.Cells.Item(8, 1).Value = "Emt"
$ws3.Cells.Item(8, 2).Value = 1
$ws3.Cells.Item(8, 3).Value = 805
$ws3.Cells.Item(8, 7).Value = 1
$ws3.Cells.Item(8, 7).HorizontalAlignment = -4108

# Row 9: Ender-IO
$ws3.Cells.Item(9, 1).Value = "Ender-IO"
$ws3.Cells.Item(9, 2).Value = 1
$ws3.Cells.Item(9, 3).Value = 1
$ws3.Cells.Item(9, 5).Value = 1
$ws3.Cells.Item(9, 5).HorizontalAlignment = -4108

# Row 10: Extra-Bees
$ws3.Cells.Item(10, 1).Value = "Extra-Bees"
$ws3.Cells.Item(10, 2).Value = 1
$ws3.Cells.Item(10, 3).Value = 41
$ws3.Cells.Item(10, 7).Value = 1
$ws3.Cells.Item(10, 7).HorizontalAlignment = -4108

# Row 11: Forbidden-Magic-01-Wands
$ws3.Cells.Item(11, 1).Value = "Forbidden-Magic-01-Wands"
$ws3.Cells.Item(11, 2).Value = 1
$ws3.Cells.Item(11, 3).Value = 132
$ws3.Cells.Item(11, 7).Value = 1
$ws3.Cells.Item(11, 7).HorizontalAlignment = -4108

# Row 12: Forestry-Frames
$ws3.Cells.Item(12, 1).Value = "Forestry-Frames"
$ws3.Cells.Item(12, 2).Value = 1
$ws3.Cells.Item(12, 3).Value = 154
$ws3.Cells.Item(12, 5).Value = 1
$ws3.Cells.Item(12, 5).HorizontalAlignment = -4108

# Row 13: Forestry
$ws3.Cells.Item(13, 1).Value = "Forestry"
$ws3.Cells.Item(13, 2).Value = 1
$ws3.Cells.Item(13, 3).Value = 25
$ws3.Cells.Item(13, 4).Value = 1
$ws3.Cells.Item(13, 4).HorizontalAlignment = -4108
$ws3.Cells.Item(13, 7).Value = 1
$ws3.Cells.Item(13, 7).HorizontalAlignment = -4108

# Row 14: Galaxy-Space
$ws3.Cells.Item(14, 1).Value = "Galaxy-Space"
$ws3.Cells.Item(14, 2).Value = 1
$ws3.Cells.Item(14, 3).Value = 2
$ws3.Cells.Item(14, 4).Value = 1
$ws3.Cells.Item(14, 4).HorizontalAlignment = -4108

# Row 15: GraviSuite
$ws3.Cells.Item(15, 1).Value = "GraviSuite"
$ws3.Cells.Item(15, 2).Value = 1
$ws3.Cells.Item(15, 3).Value = 7
$ws3.Cells.Item(15, 6).Value = 1
$ws3.Cells.Item(15, 6).HorizontalAlignment = -4108

# Row 16: Gregtech
$ws3.Cells.Item(16, 1).Value = "Gregtech"
$ws3.Cells.Item(16, 2).Value = 1
$ws3.Cells.Item(16, 3).Value = 69
$ws3.Cells.Item(16, 5).Value = 1
$ws3.Cells.Item(16, 5).HorizontalAlignment = -4108

# Row 17: Iron-Chests-Minecarts
$ws3.Cells.Item(17, 1).Value = "Iron-Chests-Minecarts"
$ws3.Cells.Item(17, 2).Value = 1
$ws3.Cells.Item(17, 3).Value = 1
$ws3.Cells.Item(17, 4).Value = 1
$ws3.Cells.Item(17, 4).HorizontalAlignment = -4108

# Row 18: Iron-Chests
$ws3.Cells.Item(18, 1).Value = "Iron-Chests"
$ws3.Cells.Item(18, 2).Value = 1
$ws3.Cells.Item(18, 3).Value = 3
$ws3.Cells.Item(18, 4).Value = 1
$ws3.Cells.Item(18, 4).HorizontalAlignment = -4108
$ws3.Cells.Item(18, 7).Value = 1
$ws3.Cells.Item(18, 7).HorizontalAlignment = -4108

# Row 19: Magic-Bees
$ws3.Cells.Item(19, 1).Value = "Magic-Bees"
$ws3.Cells.Item(19, 3).Value = 189
$ws3.Cells.Item(19, 7).Value = 1
$ws3.Cells.Item(19, 7).HorizontalAlignment = -4108

# Row 20: Minecraft
$ws3.Cells.Item(20, 1).Value = "Minecraft"
$ws3.Cells.Item(20, 2).Value = 1
$ws3.Cells.Item(20, 3).Value = 9
$ws3.Cells.Item(20, 4).Value = 1
$ws3.Cells.Item(20, 4).HorizontalAlignment = -4108

# Row 21: Open-Blocks
$ws3.Cells.Item(21, 1).Value = "Open-Blocks"
$ws3.Cells.Item(21, 3).Value = 10
$ws3.Cells.Item(21, 7).Value = 1
$ws3.Cells.Item(21, 7).HorizontalAlignment = -4108

# Row 22: Railcraft
$ws3.Cells.Item(22, 1).Value = "Railcraft"
$ws3.Cells.Item(22, 3).Value = 74
$ws3.Cells.Item(22, 6).Value = 1
$ws3.Cells.Item(22, 6).HorizontalAlignment = -4108

# Row 23: Renaming
$ws3.Cells.Item(23, 1).Value = "Renaming"
$ws3.Cells.Item(23, 3).Value = 14
$ws3.Cells.Item(23, 5).Value = 1
$ws3.Cells.Item(23, 5).HorizontalAlignment = -4108
$ws3.Cells.Item(23, 6).Value = 1
$ws3.Cells.Item(23, 6).HorizontalAlignment = -4108

# Row 24: Tainted-Magic-1
$ws3.Cells.Item(24, 1).Value = "Tainted-Magic-1"
$ws3.Cells.Item(24, 3).Value = 584
$ws3.Cells.Item(24, 7).Value = 1
$ws3.Cells.Item(24, 7).HorizontalAlignment = -4108

# Row 25: Thaumcraft-03-Alchemy
$ws3.Cells.Item(25, 1).Value = "Thaumcraft-03-Alchemy"
$ws3.Cells.Item(25, 3).Value = 138
$ws3.Cells.Item(25, 7).Value = 1
$ws3.Cells.Item(25, 7).HorizontalAlignment = -4108

# Row 26: Thaumic-Bases-01-Main
$ws3.Cells.Item(26, 1).Value = "Thaumic-Bases-01-Main"
$ws3.Cells.Item(26, 3).Value = 168
$ws3.Cells.Item(26, 7).Value = 1
$ws3.Cells.Item(26, 7).HorizontalAlignment = -4108

# Row 27: Thaumic-Exploration-01
$ws3.Cells.Item(27, 1).Value = "Thaumic-Exploration-01"
$ws3.Cells.Item(27, 3).Value = 435
$ws3.Cells.Item(27, 7).Value = 1
$ws3.Cells.Item(27, 7).HorizontalAlignment = -4108

# Row 28: ThaumicTinkerer-01
$ws3.Cells.Item(28, 1).Value = "ThaumicTinkerer-01"
$ws3.Cells.Item(28, 3).Value = 381
$ws3.Cells.Item(28, 7).Value = 1
$ws3.Cells.Item(28, 7).HorizontalAlignment = -4108

# Row 29: ThaumicTinkerer-02-Kami
$ws3.Cells.Item(29, 1).Value = "ThaumicTinkerer-02-Kami"
$ws3.Cells.Item(29, 3).Value = 174
$ws3.Cells.Item(29, 7).Value = 1
$ws3.Cells.Item(29, 7).HorizontalAlignment = -4108

# Row 30: Tinkers-Construct
$ws3.Cells.Item(30, 1).Value = "Tinkers-Construct"
$ws3.Cells.Item(30, 3).Value = 6
$ws3.Cells.Item(30, 5).Value = 1
$ws3.Cells.Item(30, 5).HorizontalAlignment = -4108

# Row 31: Warp-Theory
$ws3.Cells.Item(31, 1).Value = "Warp-Theory"
$ws3.Cells.Item(31, 3).Value = 41
$ws3.Cells.Item(31, 7).Value = 1
$ws3.Cells.Item(31, 7).HorizontalAlignment = -4108

# Row 32: Witchery
$ws3.Cells.Item(32, 1).Value = "Witchery"
$ws3.Cells.Item(32, 3).Value = 120
$ws3.Cells.Item(32, 7).Value = 1
$ws3.Cells.Item(32, 7).HorizontalAlignment = -4108

# Re-apply the AutoFilter over the (now taller/wider) table range.
$ws3.Range("A1:H43").AutoFilter()

# Drop the leftover manual sort-state (both the one nested in the old AutoFilter
# and the standalone worksheet one) -- the reverted sheet has neither.
$ws3.Sort.SortFields.Clear()

# Keep the workbook-level _xlnm._FilterDatabase name for this sheet in sync with
# the new AutoFilter range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.RefersTo() -like "*scripts*") {
        $n.RefersTo = "=scripts汉化状态!`$A`$1:`$H`$43"
    }
}

# Restore the selection that shipped with the reverted sheet.
$ws3.Activate()
$ws3.Range("F3").Select()

# --- Sheet "mod汉化状态": restore the scroll position of the view ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A88").Select()
$excel.ActiveWindow.ScrollRow = 88

# Re-activate the sheet that was active before (tab 3, "scripts汉化状态").
$ws3.Activate()